# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.183.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.328.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.06%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.691.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.317.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.096.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  +18.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -8.48%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0695"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.001.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.79%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.557.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  -6.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
